# Adds the six new case rows (76-81) for case numbers 21CRB01268 and
# 21TRD09437 to the bottom of the data table on Sheet1, mirroring the
# "No Jail dismissed checkbox" accounting columns (H/I) used by the
# existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each row: Case#, Name, Charge, Statute, Level, Plea, Result, NoJail(H), Days(I)
$rows = @(
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Dismissed", "", 0, "0"),
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Dismissed", "", 0, "0"),
    @("21CRB01268", "Bunner", "", "", "Minor Misdemeanor", "No Contest", "Guilty", 0, "0"),
    @("21TRD09437", "Bunner", "DUS", "4510.11", "M1", "Guilty", "Guilty", 0, "0"),
    @("21TRD09437", "Bunner", "1ST SPEED 1 YR SCHOOL >35MPHM4", "4511.21B1A", "M4", "Dismissed", "Guilty", 0, "0"),
    @("21TRD09437", "Bunner", "RECKLESS OPERATION 1ST IN 1 YR", "4511.20", "MM", "Guilty", "Guilty", 0, "0")
)

$startRow = 76
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    if ($data[0] -ne "") { $ws.Cells.Item($r, 1).Value = $data[0] }   # A - Case #
    if ($data[1] -ne "") { $ws.Cells.Item($r, 2).Value = $data[1] }   # B - Name
    if ($data[2] -ne "") { $ws.Cells.Item($r, 3).Value = $data[2] }   # C - Charge

    # D - Statute. Some statute citations (e.g. "4510.11", "4511.20") look
    # like plain numbers, so force the cell to Text first or Excel's COM
    # layer silently re-types the value as a number (and drops the
    # significant trailing zero). ClearFormats() afterwards drops the
    # Text number-format again so the cell keeps using the sheet's default
    # style, matching the rest of the (unstyled) worksheet.
    if ($data[3] -ne "") {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $data[3]
        $cell.ClearFormats()
    }

    if ($data[4] -ne "") { $ws.Cells.Item($r, 5).Value = $data[4] }   # E - Level
    if ($data[5] -ne "") { $ws.Cells.Item($r, 6).Value = $data[5] }   # F - Plea
    if ($data[6] -ne "") { $ws.Cells.Item($r, 7).Value = $data[6] }   # G - Result

    $ws.Cells.Item($r, 8).Value = $data[7]                            # H - No Jail (numeric)

    # I - Days, stored as text "0" (not numeric 0) in the source data.
    $iCell = $ws.Cells.Item($r, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = $data[8]
    $iCell.ClearFormats()
}
